# Insert a new price-report row for Ají "Americana (o)" at row 49 of the
# "Macroferia Regional de Talca" weekly sheet. All subsequent rows (old
# 49-139) shift down by one to rows 50-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 49; Excel shifts rows 49:139 down to 50:140
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44533
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = 100112021
$ws.Range("G49").Value = "Ají"
$ws.Range("H49").Value = "Americana (o)"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 150
$ws.Range("K49").Value = 17000
$ws.Range("L49").Value = 17000
$ws.Range("M49").Value = 17000
$ws.Range("N49").Value = "`$/caja 15 kilos"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 1133
$ws.Range("Q49").Value = 15
$ws.Range("R49").Value = "Hortaliza"
